$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on every cell we touch so Excel keeps these as literal
# strings (matching the source t="inlineStr" cells) instead of auto-converting
# "256.45" / "-0.47%" into numeric/percentage values.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('D4').NumberFormat = "@"
$ws.Range('E4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('E12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('E14').NumberFormat = "@"
$ws.Range('B15').NumberFormat = "@"
$ws.Range('C15').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('B16').NumberFormat = "@"
$ws.Range('C16').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('B17').NumberFormat = "@"
$ws.Range('C17').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('B18').NumberFormat = "@"
$ws.Range('C18').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('E18').NumberFormat = "@"
$ws.Range('B19').NumberFormat = "@"
$ws.Range('C19').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('E24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('E25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('E26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('E27').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('E50').NumberFormat = "@"

$ws.Range('D2').Value = '256.45'
$ws.Range('E2').Value = '-0.47%'
$ws.Range('D3').Value = '27.42'
$ws.Range('E3').Value = '-2.25%'
$ws.Range('D4').Value = '4.642'
$ws.Range('E4').Value = '-11.02%'
$ws.Range('D5').Value = '0.05891'
$ws.Range('E5').Value = '-0.24%'
$ws.Range('E6').Value = '-0.92%'
$ws.Range('D7').Value = '0.8643'
$ws.Range('E7').Value = '-0.42%'
$ws.Range('D8').Value = '0.9295'
$ws.Range('E8').Value = '-10.58%'
$ws.Range('D9').Value = '0.1406'
$ws.Range('E9').Value = '-0.34%'
$ws.Range('D10').Value = '0.03715'
$ws.Range('E10').Value = '3.63%'
$ws.Range('D11').Value = '0.07096'
$ws.Range('E11').Value = '-0.89%'
$ws.Range('D12').Value = '0.03228'
$ws.Range('E12').Value = '2.56%'
$ws.Range('D13').Value = '0.09206'
$ws.Range('E13').Value = '-0.22%'
$ws.Range('D14').Value = '0.001543'
$ws.Range('E14').Value = '-0.11%'
$ws.Range('B15').Value = 'One'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D15').Value = '0.0006061'
$ws.Range('E15').Value = '-94.26%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '0.006087'
$ws.Range('E16').Value = '4.57%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '3.516'
$ws.Range('E17').Value = '0.43%'
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').Value = '3.194'
$ws.Range('E18').Value = '-1.06%'
$ws.Range('B19').Value = 'BTSEToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D19').Value = '2.201'
$ws.Range('E19').Value = '-0.16%'
$ws.Range('D20').Value = '0.3101'
$ws.Range('E20').Value = '-0.61%'
$ws.Range('D21').Value = '0.1276'
$ws.Range('E21').Value = '-1.07%'
$ws.Range('D22').Value = '3.852'
$ws.Range('E22').Value = '9.19%'
$ws.Range('D23').Value = '0.04221'
$ws.Range('E23').Value = '0.51%'
$ws.Range('D24').Value = '0.001221'
$ws.Range('E24').Value = '-0.25%'
$ws.Range('D25').Value = '0.004280'
$ws.Range('E25').Value = '-5.91%'
$ws.Range('D26').Value = '0.0001200'
$ws.Range('E26').Value = '0.07%'
$ws.Range('D27').Value = '0.0001937'
$ws.Range('E27').Value = '31.64%'
$ws.Range('D40').Value = '0.03825'
$ws.Range('D41').Value = '0.006231'
$ws.Range('E41').Value = '13.15%'
$ws.Range('E42').Value = '-0.31%'
$ws.Range('D43').Value = '0.002200'
$ws.Range('E43').Value = '-4.28%'
$ws.Range('D44').Value = '0.01138'
$ws.Range('E44').Value = '5.59%'
$ws.Range('D45').Value = '0.00005476'
$ws.Range('E45').Value = '1.10%'
$ws.Range('E46').Value = '0.07%'
$ws.Range('D47').Value = '0.06021'
$ws.Range('E47').Value = '-29.54%'
$ws.Range('D48').Value = '0.002280'
$ws.Range('E48').Value = '6.70%'
$ws.Range('D49').Value = '0.00002100'
$ws.Range('E49').Value = '0.07%'
$ws.Range('D50').Value = '0.0002000'
$ws.Range('E50').Value = '0.07%'
